# "Reset fields" style edit: the previous placeholder row (ooredoo / SIM Card)
# is cleared out, an existing row's quantity is corrected, and several new
# product rows that were added through the form are appended.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old row 6 (ooredoo / SIM Card / Other / 6 / 0 / 2024-09-04 / 19:43:03).
# This shifts the old rows 7-9 up to become rows 6-8, which already match the
# desired final content for those rows.
$ws.Rows(6).Delete()

# Correct the quantity on row 3 (Master / Other / Other ...)
$ws.Range("D3").Value = 32

# Helper: write a date-like string ("YYYY-MM-DD") as literal text instead of
# letting it be auto-converted into a date serial number.
function Set-TextDate($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
}

# New row 9: LDNIO / Cable / Type C
$ws.Range("A9").Value = "LDNIO"
$ws.Range("B9").Value = "Cable"
$ws.Range("C9").Value = "Type C"
$ws.Range("D9").Value = 10
$ws.Range("E9").Value = 400
Set-TextDate $ws.Range("F9") "2024-09-10"
$ws.Range("G9").Value = "21:53:45"

# New row 10: OOREDOO / SIM Card / Other
$ws.Range("A10").Value = "OOREDOO"
$ws.Range("B10").Value = "SIM Card"
$ws.Range("C10").Value = "Other"
$ws.Range("D10").Value = 9
$ws.Range("E10").Value = 0
Set-TextDate $ws.Range("F10") "2024-09-10"
$ws.Range("G10").Value = "21:59:51"

# New row 11: jixsjixs / Cable / iPhone
$ws.Range("A11").Value = "jixsjixs"
$ws.Range("B11").Value = "Cable"
$ws.Range("C11").Value = "iPhone"
$ws.Range("D11").Value = 5
$ws.Range("E11").Value = 500
Set-TextDate $ws.Range("F11") "2024-09-11"
$ws.Range("G11").Value = "16:08:29"

# New row 12: ,ksx,ks, / Car Charger / Bluetooth
$ws.Range("A12").Value = ",ksx,ks,"
$ws.Range("B12").Value = "Car Charger"
$ws.Range("C12").Value = "Bluetooth"
$ws.Range("D12").Value = 100
$ws.Range("E12").Value = 100
Set-TextDate $ws.Range("F12") "2024-09-11"
$ws.Range("G12").Value = "16:11:40"

# New row 13: kxjoskxs / Cable / iPhone
$ws.Range("A13").Value = "kxjoskxs"
$ws.Range("B13").Value = "Cable"
$ws.Range("C13").Value = "iPhone"
$ws.Range("D13").Value = 10
$ws.Range("E13").Value = 10
Set-TextDate $ws.Range("F13") "2024-09-11"
$ws.Range("G13").Value = "16:15:14"
